$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.439.36'
$ws.Range('E2').Value = '  +1.11%  '
$ws.Range('D3').Value = '3.096.74'
$ws.Range('E3').Value = '  -0.24%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.96'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.04'
$ws.Range('E6').Value = '  +0.80%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '3.090.60'
$ws.Range('E8').Value = '  -0.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.527'
$ws.Range('E9').Value = '  -0.20%  '
$ws.Range('E10').Value = '  +6.67%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.63'
$ws.Range('E11').Value = '  -1.32%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.456'
$ws.Range('E12').Value = '  -2.27%  '
$ws.Range('E13').Value = '  +0.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.48'
$ws.Range('E14').Value = '  +5.89%  '
$ws.Range('E15').Value = '  -1.15%  '
$ws.Range('D16').Value = '3.611.69'
$ws.Range('E16').Value = '  -0.20%  '
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.12'
$ws.Range('E17').Value = '  -1.28%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '63.319.75'
$ws.Range('E18').Value = '  +1.01%  '
$ws.Range('D19').Value = '3.094.16'
$ws.Range('E19').Value = '  -0.27%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '460.66'
$ws.Range('E20').Value = '  -0.52%  '
$ws.Range('E21').Value = '  +1.03%  '
$ws.Range('E23').Value = '  -1.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.98'
$ws.Range('E24').Value = '  -2.99%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '81.18'
$ws.Range('E25').Value = '  -1.22%  '
$ws.Range('E26').Value = '  -1.98%  '
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.93'
$ws.Range('E28').Value = '  +8.22%  '
$ws.Range('B29').Value = 'FirstDigitalUSD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.67'
$ws.Range('E30').Value = '  -0.45%  '
$ws.Range('E31').Value = '  -1.95%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.80'
$ws.Range('E32').Value = '  -0.18%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '26.72'
$ws.Range('E33').Value = '  -0.79%  '
$ws.Range('E34').Value = '  -2.67%  '
$ws.Range('D35').Value = '0.0₃0849'
$ws.Range('E35').Value = '  +2.84%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.38'
$ws.Range('E36').Value = '  +7.41%  '
$ws.Range('E37').Value = '  -1.93%  '
$ws.Range('E38').Value = '  -0.92%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.01'
$ws.Range('E39').Value = '  -0.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '50.28'
$ws.Range('E40').Value = '  -1.20%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '436.91'
$ws.Range('E41').Value = '  +1.37%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.74'
$ws.Range('E42').Value = '  -0.59%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0368'
$ws.Range('E43').Value = '  -0.22%  '
$ws.Range('D44').Value = '2.864.59'
$ws.Range('E44').Value = '  -1.52%  '
$ws.Range('E45').Value = '  -1.60%  '
$ws.Range('E46').Value = '  -3.24%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '35.66'
$ws.Range('E47').Value = '  +1.66%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '123.86'
$ws.Range('E49').Value = '  -0.07%  '
$ws.Range('E50').Value = '  -1.06%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '24.14'
$ws.Range('E51').Value = '  -2.30%  '
